# Apply cryptocurrency price/volume updates per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain text in the source data (dot-grouped
# thousands like "41.803.40"). Force text via NumberFormat "@" before
# assigning so Excel does not auto-coerce numeric-looking values (e.g.
# "304.85", "1.00", "17.00") into real numbers, then restore the default
# "Normal" style so no stray formatting is left on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.803.40"
$ws.Range("E2").Value = "  -0.61%  "

Set-TextValue $ws.Range("D3") "2.270.28"
$ws.Range("E3").Value = "  +0.08%  "

Set-TextValue $ws.Range("D5") "304.85"
$ws.Range("E5").Value = "  +0.86%  "

Set-TextValue $ws.Range("D6") "92.93"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("E7").Value = "  -0.96%  "

$ws.Range("E8").Value = "  +0.01%  "

Set-TextValue $ws.Range("D9") "0.486"
$ws.Range("E9").Value = "  -0.33%  "

Set-TextValue $ws.Range("D10") "32.66"
$ws.Range("E10").Value = "  -0.99%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("E12").Value = "  -2.14%  "

Set-TextValue $ws.Range("D13") "6.66"
$ws.Range("E13").Value = "  -0.81%  "

Set-TextValue $ws.Range("D14") "2.621.16"
$ws.Range("E14").Value = "  +0.07%  "

Set-TextValue $ws.Range("D15") "14.34"
$ws.Range("E15").Value = "  +0.77%  "

Set-TextValue $ws.Range("D16") "2.274.39"
$ws.Range("E16").Value = "  +0.29%  "

Set-TextValue $ws.Range("D17") "0.786"
$ws.Range("E17").Value = "  +3.68%  "

Set-TextValue $ws.Range("D18") "41.757.77"
$ws.Range("E18").Value = "  -0.45%  "

Set-TextValue $ws.Range("D19") "12.97"
$ws.Range("E19").Value = "  +6.09%  "

$ws.Range("E20").Value = "  +0.03%  "

Set-TextValue $ws.Range("D22") "67.94"
$ws.Range("E22").Value = "  +0.72%  "

Set-TextValue $ws.Range("D23") "244.07"
$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("E26").Value = "  +0.07%  "

Set-TextValue $ws.Range("D27") "24.01"
$ws.Range("E27").Value = "  +0.14%  "

Set-TextValue $ws.Range("D28") "9.61"
$ws.Range("E28").Value = "  -1.04%  "

$ws.Range("E29").Value = "  -5.17%  "

Set-TextValue $ws.Range("D30") "34.82"
$ws.Range("E30").Value = "  +1.98%  "

Set-TextValue $ws.Range("D31") "159.63"
$ws.Range("E31").Value = "  +0.88%  "

Set-TextValue $ws.Range("D32") "5.33"
$ws.Range("E32").Value = "  +2.70%  "

Set-TextValue $ws.Range("D33") "1.00"
$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("E34").Value = "  -0.20%  "

Set-TextValue $ws.Range("D35") "3.03"
$ws.Range("E35").Value = "  -1.80%  "

Set-TextValue $ws.Range("D36") "17.00"
$ws.Range("E36").Value = "  +2.31%  "

$ws.Range("E37").Value = "  +0.87%  "

$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("E39").Value = "  +0.60%  "

$ws.Range("E40").Value = "  -0.27%  "

Set-TextValue $ws.Range("D41") "3.94"
$ws.Range("E41").Value = "  -0.99%  "

Set-TextValue $ws.Range("D42") "19.72"
$ws.Range("E42").Value = "  -2.02%  "

Set-TextValue $ws.Range("D43") "2.011.16"
$ws.Range("E43").Value = "  -1.91%  "

Set-TextValue $ws.Range("D44") "2.25"
$ws.Range("E44").Value = "  +12.57%  "

Set-TextValue $ws.Range("D45") "0.0282"
$ws.Range("E45").Value = "  +0.71%  "

Set-TextValue $ws.Range("D46") "10.28"
$ws.Range("E46").Value = "  +1.69%  "

$ws.Range("E47").Value = "  -0.17%  "

Set-TextValue $ws.Range("D48") "53.48"
$ws.Range("E48").Value = "  +2.77%  "

Set-TextValue $ws.Range("D49") "73.14"
$ws.Range("E49").Value = "  +3.44%  "

$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("E51").Value = "  -1.30%  "
